$d = $word.ActiveDocument

# --- 1. Move the "_GoBack" bookmark -----------------------------------
# It currently sits at the end of the "First name" bullet paragraph;
# it needs to move to the very start of the document, right after the
# Title paragraph's <w:pPr> (before its runs).
#
# A zero-length range exactly at document position 0 cannot be used
# directly as the bookmark anchor (the engine mis-attaches its closing
# marker to the following paragraph), so a temporary leading paragraph
# is inserted first; this gives the Title paragraph a non-zero start
# offset where the bookmark can be anchored cleanly, without touching
# (and thus merging) its existing runs. The leading paragraph is then
# removed again.

$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphBefore()

$titlePara = $d.Paragraphs(2)
$bmRange = $titlePara.Range.Duplicate
$bmRange.Collapse(1)                 # wdCollapseStart
$d.Bookmarks.Add("_GoBack", $bmRange)  # re-adding the same name moves it

$d.Paragraphs(1).Range.Delete()      # drop the temporary leading paragraph

# --- 2. Give the Title paragraph a first-line indent of 720 twips -----
$titlePara = $d.Paragraphs(1)
$titlePara.Range.ParagraphFormat.FirstLineIndent = 36   # 36 pt = 720 twips
